$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") values 45177 -> 45178 for rows 2 through 216
$ws.Range("C2:C216").Value = 45178
